$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values for columns C,D,E,F,G,I,J,M,O across rows 2-25
# (columns B,H,K,L,N remain 0 and are unchanged)
$columns = @("C","D","E","F","G","I","J","M","O")

$data = @{
    2 = @(3.652199498345238, 8.508238405086999, 14.56323218744964, 26.42636969429876, 3.624852916187835, 19.77888151176193, 10.04232738920865, 59.30097334940334, 21.15787771143988)
    3 = @(3.571563585323085, 8.508093432694199, 14.25333899226064, 26.95840458236166, 3.628753009019529, 20.00033415973993, 9.926439152986541, 55.86825641251506, 21.43837450337134)
    4 = @(3.521042329243335, 8.510840264752236, 14.06607982705072, 27.30145707343211, 3.631248651908727, 20.14740727707387, 9.858841240909246, 53.64128896996389, 21.62259288461614)
    5 = @(3.500223438909691, 8.512668214724659, 13.99063568513866, 27.44532905050861, 3.632291195082761, 20.21008204590119, 9.832211485452877, 52.7040214916205, 21.7006283909836)
    6 = @(3.496753219651361, 8.51301440687709, 13.97816366302223, 27.46946381843027, 3.632465856214326, 20.22065328266554, 9.827845635821102, 52.54659977517647, 21.71376377730898)
    7 = @(3.520762462646728, 8.510862054107967, 14.06505871444261, 27.30338093522963, 3.631262608357702, 20.14824149625992, 9.858478361758593, 53.62876879022962, 21.62363335835962)
    8 = @(3.624619376485157, 8.507597172123518, 14.45581863994214, 26.60638783832777, 3.626176809035309, 19.85290391831137, 10.00164448470932, 58.14224091554458, 21.25207287878998)
    9 = @(3.819320157650975, 8.523893710805078, 15.24144175180943, 25.37174406917842, 3.616997218142379, 19.36414896633758, 10.30958319456608, 66.04031886418925, 20.62098485490606)
    10 = @(3.955735956219571, 8.549970700274965, 15.82453854203393, 24.54861338187423, 3.610725892585874, 19.06349917552512, 10.5509475264158, 71.25750463883331, 20.22027201267074)
    11 = @(4.016136596388749, 8.564950294937338, 16.08993380379614, 24.19326256003719, 3.607973137072299, 18.94017631186484, 10.66371349979877, 73.50369623684723, 20.05246646693502)
    12 = @(4.038754326674719, 8.571074943662774, 16.19036431557736, 24.06153268000726, 3.606944934305383, 18.89547674775369, 10.70681465393855, 74.33601733690473, 19.9910787082395)
    13 = @(4.033894770322491, 8.569735717482645, 16.16873922230394, 24.08977552240582, 3.607165747419454, 18.90501355125035, 10.69751473479682, 74.15757377843495, 20.00420251941064)
    14 = @(4.018002565021641, 8.565445088521551, 16.09819815100667, 24.18236773016232, 3.607888262605317, 18.93645837453534, 10.66725161596109, 73.57253761336699, 20.04737239817896)
    15 = @(4.00823449441869, 8.562875959672908, 16.05497819874421, 24.23945495685137, 3.608332668359762, 18.95598173924225, 10.64876574412804, 73.21180868914328, 20.07409831876092)
    16 = @(3.951753865601054, 8.549054814564442, 15.80718940281959, 24.57222791877573, 3.610907789960211, 19.07183480592668, 10.54363535763028, 71.10815349537391, 20.23153614015363)
    17 = @(3.916668526471727, 8.541377214873126, 15.6551496475113, 24.7813234539644, 3.612513049900354, 19.14639708614338, 10.47988169717248, 69.78507224258263, 20.33187951157509)
    18 = @(3.896333437795435, 8.537254425784264, 15.5677186642715, 24.90338117742376, 3.613445786668447, 19.19054461845412, 10.44349330177771, 69.01210163235015, 20.39095266410312)
    19 = @(3.889422246091887, 8.535908747703614, 15.5381220015301, 24.94501300728628, 3.61376322120675, 19.20570684387251, 10.43122193063592, 68.74833120014016, 20.4111849138818)
    20 = @(3.92041958012723, 8.542164137930888, 15.67133331691549, 24.75887886011224, 3.612341192233449, 19.13832878996127, 10.48663949982221, 69.92715522196875, 20.32105668516236)
    21 = @(4.02267752831554, 8.566693046642303, 16.1189203217763, 24.15509349410073, 3.607675658365677, 18.92716743950694, 10.67613001415433, 73.7448722925551, 20.03463321709054)
    22 = @(4.088017810930619, 8.585361625159612, 16.41100607053933, 23.77705072953743, 3.604709173894296, 18.8008523557567, 10.80228691884702, 76.13354657872354, 19.86004650040362)
    23 = @(4.05328604284807, 8.575155244514002, 16.25518223491304, 23.97727166410868, 3.606284938240113, 18.86717662345346, 10.73475192101542, 74.86839148204287, 19.95204748414318)
    24 = @(3.918724238503985, 8.541807463404027, 15.66401673947914, 24.76902031789256, 3.61241885838519, 19.14197248265242, 10.48358347148836, 69.86295786215518, 20.32594537841885)
    25 = @(3.767733903116943, 8.517026837700117, 15.02746408922204, 25.69127441854412, 3.619396683582488, 19.48635703662525, 10.22351386344871, 64.00634909237995, 20.78092158775054)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}
